$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.07024041978417017
$ws.Range("D2").Value = 0.04822054979098311
$ws.Range("E2").Value = 0.108440918579241
$ws.Range("F2").Value = 2.419795505892246
$ws.Range("G2").Value = 1.920740223186542
$ws.Range("H2").Value = 1.654611435573614
$ws.Range("I2").Value = 1.751231532570721
$ws.Range("J2").Value = 0.194927321485082
$ws.Range("L2").Value = 0.2171153030101038
$ws.Range("M2").Value = 11.388345524606
$ws.Range("C3").Value = 0.06980490470707679
$ws.Range("D3").Value = 0.04559330352481794
$ws.Range("E3").Value = 0.1037931500262914
$ws.Range("F3").Value = 2.493054347851299
$ws.Range("G3").Value = 1.974577599390386
$ws.Range("H3").Value = 1.696613570917933
$ws.Range("I3").Value = 1.803894673724344
$ws.Range("J3").Value = 0.1871043495012827
$ws.Range("L3").Value = 0.2042958914198891
$ws.Range("M3").Value = 10.09659348287158
$ws.Range("C4").Value = 0.06956338470914858
$ws.Range("D4").Value = 0.04397223259358185
$ws.Range("E4").Value = 0.1010067281569178
$ws.Range("F4").Value = 2.541467271188658
$ws.Range("G4").Value = 2.010794597188053
$ws.Range("H4").Value = 1.724354593486552
$ws.Range("I4").Value = 1.838751606134807
$ws.Range("J4").Value = 0.1824993479487489
$ws.Range("L4").Value = 0.196566822843053
$ws.Range("M4").Value = 9.300405068397936
$ws.Range("C5").Value = 0.06947147026631484
$ws.Range("D5").Value = 0.0433096432405975
$ws.Range("E5").Value = 0.09988784056564981
$ws.Range("F5").Value = 2.562048048635475
$ws.Range("G5").Value = 2.026336333739607
$ws.Range("H5").Value = 1.736144787633506
$ws.Range("I5").Value = 1.85358205918304
$ws.Range("J5").Value = 0.1806713943441878
$ws.Range("L5").Value = 0.1934519324127706
$ws.Range("M5").Value = 8.975152973681645
$ws.Range("C6").Value = 0.06945660093354178
$ws.Range("D6").Value = 0.0431995008091377
$ws.Range("E6").Value = 0.09970304384607687
$ws.Range("F6").Value = 2.565516623179427
$ws.Range("G6").Value = 2.028963978679215
$ws.Range("H6").Value = 1.738131711724918
$ws.Range("I6").Value = 1.856082216790647
$ws.Range("J6").Value = 0.1803707656346774
$ws.Range("L6").Value = 0.1929367804674058
$ws.Range("M6").Value = 8.921095993921824
$ws.Range("C7").Value = 0.06956211877295715
$ws.Range("D7").Value = 0.04396330472994237
$ws.Range("E7").Value = 0.1009915716379624
$ws.Range("F7").Value = 2.541741394723658
$ws.Range("G7").Value = 2.011001043478643
$ws.Range("H7").Value = 1.724511641240269
$ws.Range("I7").Value = 1.838949090889358
$ws.Range("J7").Value = 0.1824745002533632
$ws.Range("L7").Value = 0.1965246747837313
$ws.Range("M7").Value = 9.296021879399746
$ws.Range("C8").Value = 0.0700848779988803
$ws.Range("D8").Value = 0.04731632227085214
$ws.Range("E8").Value = 0.1068241859352028
$ws.Range("F8").Value = 2.4443363706065
$ws.Range("G8").Value = 1.938640250607619
$ws.Range("H8").Value = 1.668685570032949
$ws.Range("I8").Value = 1.768861541865334
$ws.Range("J8").Value = 0.1921880159471812
$ws.Range("L8").Value = 0.2126650665753971
$ws.Range("M8").Value = 10.94356064081774
$ws.Range("C9").Value = 0.07131586659861711
$ws.Range("D9").Value = 0.05382894735670618
$ws.Range("E9").Value = 0.1188131637182295
$ws.Range("F9").Value = 2.281042671463581
$ws.Range("G9").Value = 1.822352149674913
$ws.Range("H9").Value = 1.574931304572857
$ws.Range("I9").Value = 1.651793249957294
$ws.Range("J9").Value = 0.2128715028224661
$ws.Range("L9").Value = 0.2454932644923247
$ws.Range("M9").Value = 14.15203809824186
$ws.Range("C10").Value = 0.07234686842951987
$ws.Range("D10").Value = 0.05857673014875786
$ws.Range("E10").Value = 0.1279844345268231
$ws.Range("F10").Value = 2.17864700564833
$ws.Range("G10").Value = 1.753272970291363
$ws.Range("H10").Value = 1.515961328540385
$ws.Range("I10").Value = 1.578712216401996
$ws.Range("J10").Value = 0.2291577237777318
$ws.Range("L10").Value = 0.2704068706282499
$ws.Range("M10").Value = 16.49877867328553
$ws.Range("C11").Value = 0.07284372458653365
$ws.Range("D11").Value = 0.06072892700733235
$ws.Range("E11").Value = 0.1322421195955883
$ws.Range("F11").Value = 2.136031708863143
$ws.Range("G11").Value = 1.725564354588727
$ws.Range("H11").Value = 1.491359126987845
$ws.Range("I11").Value = 1.548386120961503
$ws.Range("J11").Value = 0.2368256709060148
$ws.Range("L11").Value = 0.2819321220697759
$ws.Range("M11").Value = 17.56486815602557
$ws.Range("C12").Value = 0.07303590622541378
$ws.Range("D12").Value = 0.06154283960572116
$ws.Range("E12").Value = 0.1338672847879749
$ws.Range("F12").Value = 2.120477496070336
$ws.Range("G12").Value = 1.715620281085108
$ws.Range("H12").Value = 1.482368893791858
$ws.Range("I12").Value = 1.537331829904588
$ws.Range("J12").Value = 0.2397685451293796
$ws.Range("L12").Value = 0.2863256735764139
$ws.Range("M12").Value = 17.96842597258734
$ws.Range("C13").Value = 0.0729943365968353
$ws.Range("D13").Value = 0.06136759662302893
$ws.Range("E13").Value = 0.1335166945259871
$ws.Range("F13").Value = 2.123801220267538
$ws.Range("G13").Value = 1.717737286192573
$ws.Range("H13").Value = 1.484290495781948
$ws.Range("I13").Value = 1.539693306093604
$ws.Range("J13").Value = 0.2391329699149196
$ws.Range("L13").Value = 0.2853781189011784
$ws.Range("M13").Value = 17.88151807147693
$ws.Range("C14").Value = 0.07285945447843289
$ws.Range("D14").Value = 0.06079590958786696
$ws.Range("E14").Value = 0.1323755615801119
$ws.Range("F14").Value = 2.134740288399939
$ws.Range("G14").Value = 1.724735172434862
$ws.Range("H14").Value = 1.490612922920661
$ws.Range("I14").Value = 1.547468011954457
$ws.Range("J14").Value = 0.2370669861826542
$ws.Range("L14").Value = 0.282292987686219
$ws.Range("M14").Value = 17.59807167350306
$ws.Range("C15").Value = 0.07277736142671642
$ws.Range("D15").Value = 0.06044559492200108
$ws.Range("E15").Value = 0.1316782783753396
$ws.Range("F15").Value = 2.141517138168126
$ws.Range("G15").Value = 1.729093450061441
$ws.Range("H15").Value = 1.494528249288976
$ws.Range("I15").Value = 1.552286478785248
$ws.Range("J15").Value = 0.2358066737126023
$ws.Range("L15").Value = 0.2804071061682691
$ws.Range("M15").Value = 17.42443534760258
$ws.Range("C16").Value = 0.07231496041089258
$ws.Range("D16").Value = 0.05843592735433845
$ws.Range("E16").Value = 0.1277079546624265
$ws.Range("F16").Value = 2.181512914638731
$ws.Range("G16").Value = 1.75515977484531
$ws.Range("H16").Value = 1.517614422170169
$ws.Range("I16").Value = 1.580753665904169
$ws.Range("J16").Value = 0.228661982772195
$ws.Range("L16").Value = 0.2696576636900545
$ws.Range("M16").Value = 16.42908312237648
$ws.Range("C17").Value = 0.07203844369740864
$ws.Range("D17").Value = 0.05720112545439804
$ws.Range("E17").Value = 0.1252946053391426
$ws.Range("F17").Value = 2.207073783517359
$ws.Range("G17").Value = 1.772112355230007
$ws.Range("H17").Value = 1.532351005974164
$ws.Range("I17").Value = 1.598971851958659
$ws.Range("J17").Value = 0.2243466695177432
$ws.Range("L17").Value = 0.2631134809532369
$ws.Range("M17").Value = 15.81813657404467
$ws.Range("C18").Value = 0.07188201780876113
$ws.Range("D18").Value = 0.05649018423622465
$ws.Range("E18").Value = 0.1239145343017043
$ws.Range("F18").Value = 2.222148120223466
$ws.Range("G18").Value = 1.782211910612489
$ws.Range("H18").Value = 1.541036000098003
$ws.Range("I18").Value = 1.609724602997957
$ws.Range("J18").Value = 0.2218888972886646
$ws.Range("L18").Value = 0.2593674030674578
$ws.Range("M18").Value = 15.46659530906624
$ws.Range("C19").Value = 0.07182950375535313
$ws.Range("D19").Value = 0.05624934844237828
$ws.Range("E19").Value = 0.1234486276955522
$ws.Range("F19").Value = 2.227315597845219
$ws.Range("G19").Value = 1.785690939033671
$ws.Range("H19").Value = 1.544012278109307
$ws.Range("I19").Value = 1.613412087208019
$ws.Range("J19").Value = 0.2210608520898063
$ws.Range("L19").Value = 0.2581020840478629
$ws.Range("M19").Value = 15.34754346754295
$ws.Range("C20").Value = 0.07206760810292678
$ws.Range("D20").Value = 0.05733264621844114
$ws.Range("E20").Value = 0.1255506761624474
$ws.Range("F20").Value = 2.204314149419872
$ws.Range("G20").Value = 1.770271514438718
$ws.Range("H20").Value = 1.530760605586522
$ws.Range("I20").Value = 1.597004053006692
$ws.Range("J20").Value = 0.2248035154673858
$ws.Range("L20").Value = 0.2638082501319019
$ws.Range("M20").Value = 15.88318703547725
$ws.Range("C21").Value = 0.0728989629210588
$ws.Range("D21").Value = 0.06096385704093166
$ws.Range("E21").Value = 0.1327103859680534
$ws.Range("F21").Value = 2.131511282751177
$ws.Range("G21").Value = 1.722664719934897
$ws.Range("H21").Value = 1.488746971514445
$ws.Range("I21").Value = 1.545172655142778
$ws.Range("J21").Value = 0.2376727363076299
$ws.Range("L21").Value = 0.2831983600790409
$ws.Range("M21").Value = 17.68133022210179
$ws.Range("C22").Value = 0.073465825409329
$ws.Range("D22").Value = 0.06333079135985997
$ws.Range("E22").Value = 0.1374649900453946
$ws.Range("F22").Value = 2.087336514031662
$ws.Range("G22").Value = 1.694755348761163
$ws.Range("H22").Value = 1.463192380405019
$ws.Range("I22").Value = 1.513806293378245
$ws.Range("J22").Value = 0.2463129114210432
$ws.Range("L22").Value = 0.2960418748111522
$ws.Range("M22").Value = 18.85569600099393
$ws.Range("C23").Value = 0.07316111711440954
$ws.Range("D23").Value = 0.06206808209439885
$ws.Range("E23").Value = 0.1349202806042271
$ws.Range("F23").Value = 2.110597335634907
$ws.Range("G23").Value = 1.709353093474334
$ws.Range("H23").Value = 1.476655013519746
$ws.Range("I23").Value = 1.530314274323246
$ws.Range("J23").Value = 0.2416798360813743
$ws.Range("L23").Value = 0.2891708508901161
$ws.Range("M23").Value = 18.22896836640194
$ws.Range("C24").Value = 0.07205441494308218
$ws.Range("D24").Value = 0.05727318891288746
$ws.Range("E24").Value = 0.1254348835492962
$ws.Range("F24").Value = 2.205560600319188
$ws.Range("G24").Value = 1.771102659211181
$ws.Range("H24").Value = 1.531478963455186
$ws.Range("I24").Value = 1.597892827025241
$ws.Range("J24").Value = 0.2245969034216841
$ws.Range("L24").Value = 0.2634940944590625
$ws.Range("M24").Value = 15.85377866077937
$ws.Range("C25").Value = 0.07096071567499251
$ws.Range("D25").Value = 0.05207371149162299
$ws.Range("E25").Value = 0.1155081718221638
$ws.Range("F25").Value = 2.322181369969158
$ws.Range("G25").Value = 1.850993408100777
$ws.Range("H25").Value = 1.598577904156159
$ws.Range("I25").Value = 1.681230327683508
$ws.Range("J25").Value = 0.1806713943441878
$ws.Range("L25").Value = 0.2364787631903056
$ws.Range("M25").Value = 13.28615121475838
